# Add a new question row (row 42) to the "questions" sheet, mirroring the
# existing layout: question text (A), answer text (B), case-sensitive flag (C).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("questions")

$newRow = 42

$ws.Cells.Item($newRow, 1).Value = "From what year World War II ends?"

# Force column B to text so the numeric-looking answer "1945" is stored as a
# string (matching the sheet's existing answer column, which is all text),
# not auto-coerced to a number. Reset the style afterwards so no extra
# number-format is left applied to the cell.
$ws.Cells.Item($newRow, 2).NumberFormat = "@"
$ws.Cells.Item($newRow, 2).Value = "1945"
$ws.Cells.Item($newRow, 2).Style = "Normal"

$ws.Cells.Item($newRow, 3).Value = $true
